$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down.
$ws.Rows.Item(1).Insert()

# Populate new row 1 with the BY4743 normalization data.
$ws.Range("A1").Value = "BY4743"
$ws.Range("C1").Value = 0.08
$ws.Range("D1").Value = "(0.15)"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").Value = 0.05
$ws.Range("F1").Value = "(0.09)"
$ws.Range("F1").NumberFormat = "@"
$ws.Range("G1").Value = 0.06
$ws.Range("H1").Value = "(0.32)"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("I1").Value = 0.03
$ws.Range("J1").Value = "(0.31)"
$ws.Range("J1").NumberFormat = "@"
$ws.Range("K1").Value = 0.22

# Match selection noted in diff.
$ws.Range("B14").Select()
